$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.609.21'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '3.843.69'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'516.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.46%  '
$ws.Range("D6").Value = "'140.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.65%  '
$ws.Range("E7").Value = '  -3.12%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("E10").Value = '  -5.23%  '
$ws.Range("E11").Value = '  -8.89%  '
$ws.Range("D12").Value = "'41.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.36%  '
$ws.Range("D13").Value = "'10.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").Value = '4.456.24'
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = "'21.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.96%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.837.03'
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").Value = "'13.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").Value = '68.498.83'
$ws.Range("E20").Value = '  -1.16%  '
$ws.Range("D21").Value = "'413.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.88%  '
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").Value = "'12.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").Value = "'13.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.62%  '
$ws.Range("D25").Value = "'86.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.61%  '
$ws.Range("E26").Value = '  +4.40%  '
$ws.Range("E27").Value = '  -7.06%  '
$ws.Range("D28").Value = "'35.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.16%  '
$ws.Range("D29").Value = "'13.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("D30").Value = "'676.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("D31").Value = "'6.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +13.47%  '
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("D33").Value = "'0.123"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.49%  '
$ws.Range("D34").Value = "'66.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.95%  '
$ws.Range("D35").Value = "'0.440"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.40%  '
$ws.Range("D36").Value = '0.0₃0844'
$ws.Range("E36").Value = '  -6.48%  '
$ws.Range("D37").Value = "'39.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("D38").Value = "'3.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.85%  '
$ws.Range("E39").Value = '  -3.62%  '
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").Value = "'0.0471"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = '  +4.92%  '
$ws.Range("D44").Value = "'2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("D45").Value = "'3.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("D47").Value = "'0.000283"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +17.34%  '
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").Value = "'142.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("D51").Value = "'8.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.67%  '
